$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "JULIA" -> "MARIA JULIA" for Julia Bonilla Gomez (row 7, Nombre column C)
$ws.Range("C7").Value = "MARIA JULIA"

# Clear/apply "No Fill" formatting on D7 (last touched cell), moving the selection there
$ws.Range("D7").Interior.ColorIndex = -4142
$ws.Range("D7").Select()
